# Emma.xlsx re-upload edit
# - "emma" sheet becomes the selected/active tab (was "Formatted")
# - columns A and B on "emma" are un-hidden and given new widths
# - the RAND() driven column A values are left to recalc naturally
# - the "Formatted" sheet loses its tabSelected flag

$wb = $excel.ActiveWorkbook

$emma = $wb.Worksheets.Item("emma")
$formatted = $wb.Worksheets.Item("Formatted")

# Un-hide columns A (1) and B (2) on "emma" and resize them.
$emma.Columns.Item(1).Hidden = $false
$emma.Columns.Item(2).Hidden = $false
$emma.Columns.Item(1).ColumnWidth = 12.3
$emma.Columns.Item(2).ColumnWidth = 12.3

# Make "emma" the active/selected sheet (was "Formatted"), which flips
# tabSelected on both sheets and updates the workbook's active-tab index.
$emma.Activate()
